# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets.
# Row -> (old, new) mapping, identical update applied to both sheets.
$updates = @{
    3  = 1006
    4  = 240
    5  = 1367
    6  = 8417
    7  = 52
    10 = 235
    12 = 3381
    13 = 43
    14 = 337
    16 = 905
    17 = 138
    18 = 1089
    19 = 291
    20 = 147
    21 = 1972
}

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
